# Add 2022-Q4 data.
#
# Before: 总计, 2022-Q3
# After : 总计, 2022-Q4, 2022-Q3
#
# The existing "2022-Q3" sheet is first duplicated (the duplicate is
# placed right after it and keeps the original look/figures, becoming
# the "new" 2022-Q3 tab); the original sheet is then updated in place
# with the Q4 figures and restyled to match the "总计" sheet's header
# look, and renamed to "2022-Q4".

$wb    = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)   # "总计"
$q3    = $wb.Worksheets.Item(2)   # "2022-Q3"

# --- 1. Duplicate "2022-Q3" right after itself; this duplicate keeps
#        the old figures/format and becomes the "2022-Q3" tab. --------
$q3.Copy($null, $q3)
$q3Dup = $wb.Worksheets.Item(3)

# --- 2. Turn the original sheet into "2022-Q4" (rename the original
#        first so the duplicate can take over the "2022-Q3" name). ----
$q3.Name = "2022-Q4"
$q3Dup.Name = "2022-Q3"

$q3.Range("D2").Value = "'4.76"
$q3.Range("E2").Value = "'92.90"
$q3.Range("F2").Value = "'3.34"
$q3.Range("G2").Value = "'0.1590"
$q3.Range("H2").Value = 9

# --- 3. Match "总计"'s header/border style and page margins. ----------
$total.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)

$q3.PageSetup.LeftMargin   = 54
$q3.PageSetup.RightMargin  = 54
$q3.PageSetup.TopMargin    = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# --- 4. Update the "总计" summary sheet: rename the Q3 row to Q4 with
#        its new total, and append a new row for the Q3 total. --------
$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 0.16

$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.18
